$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for data rows 2 through 12
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13)
$ws.Range("C2:C12").Value = 45212
